$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was inserted at the top of the data (row 9),
# pushing the previously-existing rows 9-79 down to rows 10-80.
$ws.Rows("9:9").Insert()

# Populate the newly-inserted row 9 with the new record's values.
$ws.Range("A9").Value = 9
$ws.Range("B9").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C9").Value = "Metropolitana"
$ws.Range("D9").Value = 45069
$ws.Range("E9").Value = 13
$ws.Range("F9").Value = "Fruta"
$ws.Range("G9").Value = 100104
$ws.Range("H9").Value = "Frutos de pepita"
$ws.Range("I9").Value = 100104003
$ws.Range("J9").Value = "Membrillo"
$ws.Range("K9").Value = "Champion"
$ws.Range("L9").Value = "Primera"
$ws.Range("M9").Value = 380
$ws.Range("N9").Value = 9000
$ws.Range("O9").Value = 9000
$ws.Range("P9").Value = 9000
$ws.Range("Q9").Value = "`$/caja 15 kilos empedrada"
$ws.Range("R9").Value = "Provincia de Curicó"
$ws.Range("S9").Value = 600
$ws.Range("T9").Value = 15
